$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect repulled data / mean calculation
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -1
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = -1
